# Generate Report for Handoff
# Rotate the report's file identifiers (old GUID -> new GUID) and refresh
# the handoff/target timestamps, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuidFile   = "99d98153-0281-4b48-b8b6-63859d6a83d6"
$newGuidFile   = "ae71c06a-3c8b-4a95-b378-17294e303fb0"

$oldZhCnHash   = "f86a4668050a9217bd94b4cced0c2072aefcf175"
$newZhCnHash   = "b0648dce9773cc3393006c592012f66909d3cb13"
$oldDeDeHash   = "f86a4668050a9217bd94b4cced0c2072aefcf175"
$newDeDeHash   = "b0648dce9773cc3393006c592012f66909d3cb13"

$newMdName     = $newGuidFile + ".md"
$newMdPath     = "e2e\" + $newGuidFile + ".md"
$newZhCnXlf    = $newGuidFile + "." + $newZhCnHash + ".zh-cn.xlf"
$newDeDeXlf    = $newGuidFile + "." + $newDeDeHash + ".de-de.xlf"

$newHoDate       = "2016-08-22 06:56:11"
$newZhCnHoDate   = "2016-08-22 06:56:02"
$newDeDeHoDate   = "2016-08-22 06:56:11"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newMdPath
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdPath
}
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMdName
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
$wsZhCn.Range("G2").Value = $newZhCnXlf
$wsZhCn.Range("H2").Value = $newZhCnHoDate

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMdName
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}
$wsDeDe.Range("G2").Value = $newDeDeXlf
$wsDeDe.Range("H2").Value = $newDeDeHoDate

Write-Output "Report identifiers and timestamps updated."
